$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.034.59'
$ws.Range('D3').Value = '1.563.01'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.29'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.05'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.786.27'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.566.16'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '27.039.54'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.95'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '0.0₃0709'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.92'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.38'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.40'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.60'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.06'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  +2.70%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('E33').Value = '  +2.72%  '
$ws.Range('D34').Value = '1.429.08'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('E36').Value = '  +7.91%  '
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.533'
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.90'
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.810'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.71'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').Value = '1.700.69'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.03'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  +0.33%  '
